$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "D" (Price) column cells we touch keep their original plain-text
# representation (e.g. "1.001", "317.83") instead of being auto-converted to
# numbers/dates by Excel's smart input parsing.
$textCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values scraped for this run.
$ws.Range("D2").Value = '25.111.43'
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").Value = '1.708.94'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '317.83'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.4012'
$ws.Range("E7").Value = '  +2.67%  '
$ws.Range("D8").Value = '0.4042'
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Value = '1.473'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").Value = '52.90'
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("D11").Value = '1.004'
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = '0.08830'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '26.06'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = '7.502'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.00001357'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '7.983'
$ws.Range("E16").Value = '  -4.05%  '
$ws.Range("D17").Value = '1.705.64'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '96.32'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = '0.07205'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '20.85'
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").Value = '7.315'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '25.098.13'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").Value = '2.399'
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("D26").Value = '2.953'
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("D27").Value = '23.61'
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").Value = '6.205'
$ws.Range("E28").Value = '  +15.42%  '
$ws.Range("D29").Value = '163.20'
$ws.Range("E29").Value = '  -3.14%  '
$ws.Range("D30").Value = '152.17'
$ws.Range("E30").Value = '  +3.24%  '
$ws.Range("D31").Value = '8.351'
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").Value = '2.695'
$ws.Range("E32").Value = '  +23.65%  '
$ws.Range("D33").Value = '1.893.12'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '0.08622'
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").Value = '0.03179'
$ws.Range("E35").Value = '  +3.26%  '
$ws.Range("D36").Value = '1.048'
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("D37").Value = '7.216'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("D38").Value = '0.2919'
$ws.Range("E38").Value = '  +3.88%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '11.14'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.09783'
$ws.Range("E40").Value = '  +6.43%  '
$ws.Range("D41").Value = '0.8311'
$ws.Range("E41").Value = '  +3.81%  '
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").Value = '1.482'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").Value = '17.16'
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").Value = '2.690'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '0.7403'
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '0.09182'
$ws.Range("E47").Value = '  +11.87%  '
$ws.Range("B48").Value = 'Flow'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D48").Value = '1.441'
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '4.254'
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '140.49'
$ws.Range("E51").Value = '  -0.12%  '
